$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh re-sorted the price rows by date; row contents (D, J, K, L, M, P)
# were shuffled accordingly while A/B/C/E/F/G/H/I/N/O/Q/R stay constant per row.

$ws.Cells.Item(2, 4).Value = 45203   # D2
$ws.Cells.Item(2, 10).Value = 800  # J2
$ws.Cells.Item(2, 11).Value = 1800  # K2
$ws.Cells.Item(2, 12).Value = 2000  # L2
$ws.Cells.Item(2, 13).Value = 1900  # M2
$ws.Cells.Item(2, 16).Value = 1900  # P2

$ws.Cells.Item(3, 4).Value = 44883   # D3
$ws.Cells.Item(3, 10).Value = 290  # J3
$ws.Cells.Item(3, 11).Value = 1400  # K3
$ws.Cells.Item(3, 12).Value = 1500  # L3
$ws.Cells.Item(3, 13).Value = 1434  # M3
$ws.Cells.Item(3, 16).Value = 1434  # P3

$ws.Cells.Item(4, 4).Value = 45204   # D4
$ws.Cells.Item(4, 10).Value = 1200  # J4
$ws.Cells.Item(4, 11).Value = 1600  # K4
$ws.Cells.Item(4, 12).Value = 1700  # L4
$ws.Cells.Item(4, 13).Value = 1650  # M4
$ws.Cells.Item(4, 16).Value = 1650  # P4

$ws.Cells.Item(5, 4).Value = 44537   # D5
$ws.Cells.Item(5, 10).Value = 800  # J5
$ws.Cells.Item(5, 11).Value = 1300  # K5
$ws.Cells.Item(5, 12).Value = 1400  # L5
$ws.Cells.Item(5, 13).Value = 1350  # M5
$ws.Cells.Item(5, 16).Value = 1350  # P5

$ws.Cells.Item(6, 4).Value = 44200   # D6
$ws.Cells.Item(6, 10).Value = 1500  # J6
$ws.Cells.Item(6, 11).Value = 1400  # K6
$ws.Cells.Item(6, 12).Value = 1500  # L6
$ws.Cells.Item(6, 13).Value = 1450  # M6
$ws.Cells.Item(6, 16).Value = 1450  # P6

$ws.Cells.Item(7, 4).Value = 45210   # D7
$ws.Cells.Item(7, 10).Value = 550  # J7
$ws.Cells.Item(7, 11).Value = 1500  # K7
$ws.Cells.Item(7, 12).Value = 1600  # L7
$ws.Cells.Item(7, 13).Value = 1536  # M7
$ws.Cells.Item(7, 16).Value = 1536  # P7

$ws.Cells.Item(8, 4).Value = 44893   # D8
$ws.Cells.Item(8, 10).Value = 3300  # J8
$ws.Cells.Item(8, 11).Value = 1200  # K8
$ws.Cells.Item(8, 12).Value = 1300  # L8
$ws.Cells.Item(8, 13).Value = 1261  # M8
$ws.Cells.Item(8, 16).Value = 1261  # P8

$ws.Cells.Item(9, 4).Value = 44210   # D9
$ws.Cells.Item(9, 10).Value = 1450  # J9
$ws.Cells.Item(9, 11).Value = 1600  # K9
$ws.Cells.Item(9, 12).Value = 1700  # L9
$ws.Cells.Item(9, 13).Value = 1650  # M9
$ws.Cells.Item(9, 16).Value = 1650  # P9

$ws.Cells.Item(10, 4).Value = 44907   # D10
$ws.Cells.Item(10, 10).Value = 2300  # J10
$ws.Cells.Item(10, 11).Value = 900  # K10
$ws.Cells.Item(10, 12).Value = 1000  # L10
$ws.Cells.Item(10, 13).Value = 952  # M10
$ws.Cells.Item(10, 16).Value = 952  # P10

$ws.Cells.Item(11, 4).Value = 45132   # D11
$ws.Cells.Item(11, 10).Value = 170  # J11
$ws.Cells.Item(11, 11).Value = 2200  # K11
$ws.Cells.Item(11, 12).Value = 2500  # L11
$ws.Cells.Item(11, 13).Value = 2359  # M11
$ws.Cells.Item(11, 16).Value = 2359  # P11

$ws.Cells.Item(12, 4).Value = 44638   # D12
$ws.Cells.Item(12, 10).Value = 800  # J12
$ws.Cells.Item(12, 11).Value = 2500  # K12
$ws.Cells.Item(12, 12).Value = 2800  # L12
$ws.Cells.Item(12, 13).Value = 2650  # M12
$ws.Cells.Item(12, 16).Value = 2650  # P12

$ws.Cells.Item(13, 4).Value = 45205   # D13
$ws.Cells.Item(13, 10).Value = 3500  # J13
$ws.Cells.Item(13, 11).Value = 1400  # K13
$ws.Cells.Item(13, 12).Value = 1500  # L13
$ws.Cells.Item(13, 13).Value = 1457  # M13
$ws.Cells.Item(13, 16).Value = 1457  # P13

$ws.Cells.Item(14, 4).Value = 44895   # D14
$ws.Cells.Item(14, 10).Value = 200  # J14
$ws.Cells.Item(14, 11).Value = 1200  # K14
$ws.Cells.Item(14, 12).Value = 1300  # L14
$ws.Cells.Item(14, 13).Value = 1255  # M14
$ws.Cells.Item(14, 16).Value = 1255  # P14

$ws.Cells.Item(15, 4).Value = 45062   # D15
$ws.Cells.Item(15, 10).Value = 1700  # J15
$ws.Cells.Item(15, 11).Value = 2800  # K15
$ws.Cells.Item(15, 12).Value = 3000  # L15
$ws.Cells.Item(15, 13).Value = 2900  # M15
$ws.Cells.Item(15, 16).Value = 2900  # P15

$ws.Cells.Item(16, 4).Value = 45212   # D16
$ws.Cells.Item(16, 10).Value = 750  # J16
$ws.Cells.Item(16, 11).Value = 1400  # K16
$ws.Cells.Item(16, 12).Value = 1500  # L16
$ws.Cells.Item(16, 13).Value = 1440  # M16
$ws.Cells.Item(16, 16).Value = 1440  # P16

$ws.Cells.Item(17, 4).Value = 44175   # D17
$ws.Cells.Item(17, 10).Value = 1400  # J17
$ws.Cells.Item(17, 11).Value = 1900  # K17
$ws.Cells.Item(17, 12).Value = 2000  # L17
$ws.Cells.Item(17, 13).Value = 1950  # M17
$ws.Cells.Item(17, 16).Value = 1950  # P17
